# Update "想去人数" (F column) values for two entries that appear on both
# the "展览" sheet (sheet1) and the "全部类型" sheet (sheet4).
#   - Row for 2024.03.16 "南宁·草莓动漫节": 876 -> 894
#   - Row for 2024.03.30 "南宁·第一届ANE·DACG动漫嘉年华": 532 -> 534

$wb = $excel.ActiveWorkbook

# Sheet "展览": F4 is the 南宁·草莓动漫节 row, F5 is the DACG动漫嘉年华 row.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 894
$wsExhibit.Range("F5").Value = 534

# Sheet "全部类型": F4 is the 南宁·草莓动漫节 row, F6 is the DACG动漫嘉年华 row.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 894
$wsAll.Range("F6").Value = 534
